$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row needs to be inserted right before the
# existing row 313 (the data set is ordered, most-recent-first groups of
# "Primera"/"Segunda" quality rows per reporting date). Inserting the row
# pushes the old rows 313-326 down to 314-327 and updates the used range.
$ws.Rows("313:313").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A313").Value2 = 7
$ws.Range("B313").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C313").Value2 = "Ñuble"
$ws.Range("D313").Value2 = 45075
$ws.Range("E313").Value2 = 16
$ws.Range("F313").Value2 = "Fruta"
$ws.Range("G313").Value2 = 100108
$ws.Range("H313").Value2 = "Tropicales y subtropicales"
$ws.Range("I313").Value2 = 100108005
$ws.Range("J313").Value2 = "Piña"
$ws.Range("K313").Value2 = "Caramelo"
$ws.Range("L313").Value2 = "Primera"
$ws.Range("M313").Value2 = 30
$ws.Range("N313").Value2 = 15000
$ws.Range("O313").Value2 = 15000
$ws.Range("P313").Value2 = 15000
$ws.Range("Q313").Value2 = "$/caja 12 unidades"
$ws.Range("R313").Value2 = "Ecuador"
$ws.Range("S313").Value2 = 1250
$ws.Range("T313").Value2 = 12
